$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# More specs have been written: cells that were marked "X" (not yet covered)
# are now marked "Y-Test" (spec written) for rows 14-17, and a new spec
# cell is added at F16.
$ws.Range("F14").Value = "Y-Test"
$ws.Range("C15").Value = "Y-Test"
$ws.Range("F15").Value = "Y-Test"
$ws.Range("C16").Value = "Y-Test"
$ws.Range("F16").Value = "Y-Test"
$ws.Range("F16").HorizontalAlignment = -4108
$ws.Range("C17").Value = "Y-Test"
$ws.Range("F17").Value = "Y-Test"

$ws.Range("I18").Select()
